# Region 2.xlsx edit
#  - Rename sheet "Iligan" -> "Ilagan" (spelling correction for the city name)
#  - Make the "Ilagan" sheet the active tab (was "Tuguegarao")
#  - Move the selection on the "Ilagan" sheet from E35 to D36

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Iligan")
$ws.Name = "Ilagan"

$ws.Activate()
$ws.Range("D36").Select()

Write-Output "renamed sheet, activated Ilagan, selected D36"
